$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I: "Маг-майнер" -> "Тип"
$ws.Cells.Item(1, 9).Value = "Тип"

# Column I (2..5) becomes text values "МАЙНОР1" / "МАЙНОР3" instead of numbers 1/2
$ws.Cells.Item(3, 9).Value = "МАЙНОР3"
$ws.Cells.Item(2, 9).Value = "МАЙНОР1"
$ws.Cells.Item(5, 9).Value = "МАЙНОР3"
$ws.Cells.Item(4, 9).Value = "МАЙНОР1"

# Column I formatting: new font (Helvetica Neue, 14pt, dark grey) for the value cells
$valueRange = $ws.Range("I2:I5")
$valueRange.Font.Name = "Helvetica Neue"
$valueRange.Font.Size = 14
$valueRange.Font.Color = 3355443

# Column width for I and row heights for rows 4 and 5
$ws.Columns.Item(9).ColumnWidth = 18.83
$ws.Rows.Item(4).RowHeight = 18
$ws.Rows.Item(5).RowHeight = 18

# New selection reflects where the author left off working
$ws.Range("F8").Select()
